$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 0.061724
    "H2" = 0.185172
    "I2" = 0.09652262708432048
    "J2" = 0.09652262708432047
    "M2" = 6.977989333333333
    "N2" = 20.933968
    "O2" = 0.08060734116444308
    "P2" = 0.0806073411644431
    "Q2" = 0.4307094136106667
    "R2" = 3.876384722496
    "S2" = 0.007780432331474135
    "T2" = 0.007780432331474136

    "G3" = 0.061724
    "H3" = 0.185172
    "I3" = 0.09652262708432048
    "J3" = 0.09652262708432047
    "O3" = 0.8031574845597853
    "P3" = 0.8031574845597854
    "Q3" = 4.291513455406666
    "R3" = 38.62362109866
    "S3" = 0.07752287037214504
    "T3" = 0.07752287037214504

    "G4" = 0.061724
    "H4" = 0.185172
    "I4" = 0.09652262708432048
    "J4" = 0.09652262708432047
    "O4" = 0.1162351742757715
    "P4" = 0.1162351742757716
    "Q4" = 0.6210797060173334
    "R4" = 5.589717354156
    "S4" = 0.0112193243807013
    "T4" = 0.0112193243807013

    "I5" = 0.8735221647273214
    "J5" = 0.8735221647273215
    "M5" = 6.977989333333333
    "N5" = 20.933968
    "O5" = 0.08060734116444308
    "P5" = 0.0806073411644431
    "Q5" = 3.897886233628444
    "R5" = 35.080976102656
    "S5" = 0.07041229914687805
    "T5" = 0.07041229914687806

    "I6" = 0.8735221647273214
    "J6" = 0.8735221647273215
    "O6" = 0.8031574845597853
    "P6" = 0.8031574845597854
    "S6" = 0.7015758645296138
    "T6" = 0.7015758645296141

    "I7" = 0.8735221647273214
    "J7" = 0.8735221647273215
    "O7" = 0.1162351742757715
    "P7" = 0.1162351742757716
    "Q7" = 5.620722370046222
    "R7" = 50.586501330416
    "S7" = 0.1015340010508294
    "T7" = 0.1015340010508294

    "I8" = 0.02995520818835809
    "J8" = 0.02995520818835809
    "M8" = 6.977989333333333
    "N8" = 20.933968
    "O8" = 0.08060734116444308
    "P8" = 0.0806073411644431
    "Q8" = 0.1336680376728889
    "R8" = 1.203012339056
    "S8" = 0.0024146096860909
    "T8" = 0.0024146096860909

    "I9" = 0.02995520818835809
    "J9" = 0.02995520818835809
    "O9" = 0.8031574845597853
    "P9" = 0.8031574845597854
    "S9" = 0.02405874965802637
    "T9" = 0.02405874965802637

    "I10" = 0.02995520818835809
    "J10" = 0.02995520818835809
    "O10" = 0.1162351742757715
    "P10" = 0.1162351742757716
    "Q10" = 0.1927482959934445
    "S10" = 0.003481848844240821
    "T10" = 0.003481848844240822
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
